$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10: fix the date ---
$ws.Range("C10").Value = 43511

# --- Row 11: fill in 4th log entry ---
$ws.Range("C11").Value = 43511
$ws.Range("D11").Value = 0.46527777777777773
$ws.Range("E11").Value = 0.52083333333333337
$ws.Range("G11").Value = 80
$ws.Range("H11").Value = "Summarising"

$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4163)

# --- Row 12: sequence number ---
$ws.Range("B12").Value = 5

# --- sheet view changes ---
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G13").Select()

$wb.Save()
